# Applies the odds/score updates for the 2024-10-20 FlashScore weekly workbook.
# Each assignment below mirrors one cell delta from the authoritative diff
# (row numbers refer to Sheet1 rows as in the OOXML sheet data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25
$ws.Range("G25").Value = 2.25
$ws.Range("I25").Value = 3.5
$ws.Range("J25").Value = 3.1
$ws.Range("L25").Value = 4.33
$ws.Range("W25").Value = 6
$ws.Range("X25").Value = 9.5
$ws.Range("AA25").Value = 23
$ws.Range("AE25").Value = 21
$ws.Range("AH25").Value = 7.5
$ws.Range("AI25").Value = 15
$ws.Range("AM25").Value = 41
$ws.Range("AQ25").Value = 51
$ws.Range("AX25").Value = 5
$ws.Range("AY25").Value = 21

# Row 26
$ws.Range("G26").Value = 1.8
$ws.Range("H26").Value = 3.3
$ws.Range("I26").Value = 5.75
$ws.Range("J26").Value = 2.5
$ws.Range("L26").Value = 6
$ws.Range("M26").Value = 1.13
$ws.Range("N26").Value = 6
$ws.Range("O26").Value = 1.53
$ws.Range("P26").Value = 2.38
$ws.Range("X26").Value = 6.5
$ws.Range("Z26").Value = 13
$ws.Range("AA26").Value = 19
$ws.Range("AH26").Value = 10
$ws.Range("AI26").Value = 26
$ws.Range("AK26").Value = 67
$ws.Range("AM26").Value = 67
$ws.Range("AN26").Value = 3.5
$ws.Range("AO26").Value = 10
$ws.Range("AR26").Value = 67
$ws.Range("AY26").Value = 34
$ws.Range("BA26").Value = 151

# Row 27
$ws.Range("H27").Value = 3.1
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 3
$ws.Range("L27").Value = 4.75
$ws.Range("M27").Value = 1.13
$ws.Range("N27").Value = 6
$ws.Range("S27").Value = 1.62
$ws.Range("T27").Value = 2.2
$ws.Range("U27").Value = 2.25
$ws.Range("V27").Value = 1.57
$ws.Range("AC27").Value = 6
$ws.Range("AD27").Value = 6
$ws.Range("AT27").Value = 2.2
$ws.Range("BB27").Value = 151

# Row 28
$ws.Range("Z28").Value = 13
$ws.Range("AA28").Value = 19

# Row 43
$ws.Range("G43").Value = 2.35
$ws.Range("I43").Value = 2.9
$ws.Range("L43").Value = 3.4
$ws.Range("Q43").Value = 1.83
$ws.Range("R43").Value = 1.98
$ws.Range("S43").Value = 1.36
$ws.Range("T43").Value = 3
$ws.Range("U43").Value = 1.62
$ws.Range("V43").Value = 2.2
$ws.Range("W43").Value = 9
$ws.Range("Z43").Value = 23
$ws.Range("AC43").Value = 11
$ws.Range("AD43").Value = 6.5
$ws.Range("AE43").Value = 12
$ws.Range("AK43").Value = 29
$ws.Range("AT43").Value = 3
$ws.Range("AY43").Value = 15

# Row 45
$ws.Range("Q45").Value = 2.05
$ws.Range("R45").Value = 1.85

# Row 48
$ws.Range("G48").Value = 1.42
$ws.Range("H48").Value = 4.1
$ws.Range("J48").Value = 2
$ws.Range("K48").Value = 2.2
$ws.Range("M48").Value = 1.07
$ws.Range("N48").Value = 9
$ws.Range("Q48").Value = 2.2
$ws.Range("R48").Value = 1.65
$ws.Range("S48").Value = 1.44
$ws.Range("T48").Value = 2.63
$ws.Range("U48").Value = 2.5
$ws.Range("V48").Value = 1.5
$ws.Range("Y48").Value = 9
$ws.Range("Z48").Value = 8.5
$ws.Range("AH48").Value = 17
$ws.Range("AJ48").Value = 29
$ws.Range("AK48").Value = 126
$ws.Range("AL48").Value = 81
$ws.Range("AO48").Value = 7
$ws.Range("AT48").Value = 2.63

# Row 52
$ws.Range("G52").Value = 2.1
$ws.Range("H52").Value = 3.25
$ws.Range("I52").Value = 3.4
$ws.Range("J52").Value = 2.75
$ws.Range("L52").Value = 4
$ws.Range("X52").Value = 10
$ws.Range("Z52").Value = 19
$ws.Range("AD52").Value = 6
$ws.Range("AE52").Value = 13
$ws.Range("AG52").Value = 201
$ws.Range("AH52").Value = 10
$ws.Range("AI52").Value = 17
$ws.Range("AJ52").Value = 12
$ws.Range("AK52").Value = 34
$ws.Range("AM52").Value = 34
$ws.Range("AO52").Value = 12
$ws.Range("AY52").Value = 19
$ws.Range("AZ52").Value = 26
$ws.Range("BA52").Value = 51

# Row 53
$ws.Range("I53").Value = 1.9
$ws.Range("K53").Value = 2.38
$ws.Range("AC53").Value = 15
$ws.Range("AI53").Value = 11
$ws.Range("AS53").Value = 126
$ws.Range("AX53").Value = 4.33

# Row 54
$ws.Range("G54").Value = 1.67
$ws.Range("I54").Value = 5.25
$ws.Range("K54").Value = 2.2
$ws.Range("S54").Value = 1.4
$ws.Range("T54").Value = 2.75
$ws.Range("AC54").Value = 10
$ws.Range("AG54").Value = 301
$ws.Range("AI54").Value = 26
$ws.Range("AJ54").Value = 17
$ws.Range("AN54").Value = 3.6
$ws.Range("AO54").Value = 8.5
$ws.Range("AT54").Value = 2.75
$ws.Range("BA54").Value = 101

# Row 55
$ws.Range("G55").Value = 1.42
$ws.Range("N55").Value = 17

# Row 56
$ws.Range("I56").Value = 2

# Row 59
$ws.Range("G59").Value = 3.7
$ws.Range("J59").Value = 4.75
$ws.Range("AR59").Value = 151

# Row 61
$ws.Range("G61").Value = 2.25
$ws.Range("H61").Value = 2.8
$ws.Range("I61").Value = 3.75
$ws.Range("J61").Value = 3.2
$ws.Range("L61").Value = 4.5
$ws.Range("AA61").Value = 26
$ws.Range("AC61").Value = 5
$ws.Range("AE61").Value = 23
$ws.Range("AH61").Value = 7
$ws.Range("AO61").Value = 15

# Row 75
$ws.Range("Q75").Value = 2.5
$ws.Range("R75").Value = 1.5
$ws.Range("U75").Value = 2.25
$ws.Range("V75").Value = 1.57
$ws.Range("AD75").Value = 6
$ws.Range("AH75").Value = 9

# Row 120
$ws.Range("G120").Value = 2.63
$ws.Range("I120").Value = 2.9
$ws.Range("K120").Value = 1.95
$ws.Range("M120").Value = 1.1
$ws.Range("N120").Value = 7
$ws.Range("O120").Value = 1.44
$ws.Range("P120").Value = 2.63
$ws.Range("Q120").Value = 2.4
$ws.Range("R120").Value = 1.53
$ws.Range("S120").Value = 1.53
$ws.Range("T120").Value = 2.38
$ws.Range("U120").Value = 2
$ws.Range("V120").Value = 1.73
$ws.Range("W120").Value = 7
$ws.Range("AB120").Value = 41
$ws.Range("AC120").Value = 7
$ws.Range("AE120").Value = 17
$ws.Range("AH120").Value = 7.5
$ws.Range("AS120").Value = 251
$ws.Range("AT120").Value = 2.38

# Row 144
$ws.Range("B144").Value = "21/10/2024"

# Row 165
$ws.Range("M165").Value = 1.08
$ws.Range("N165").Value = 8
$ws.Range("AC165").Value = 8
$ws.Range("AG165").Value = 351
$ws.Range("BC165").Value = 251
